$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-09 Sunday" "2025-03-10 Monday"

Replace-Text "36×64=" "80×76="
Replace-Text "26×52=" "34×15="
Replace-Text "48×44=" "82×71="
Replace-Text "95×49=" "55×65="
Replace-Text "94×51=" "76×12="

Replace-Text "27×12=" "90×35="
Replace-Text "65×22=" "26×25="
Replace-Text "99×89=" "99×52="
Replace-Text "56×21=" "33×79="
Replace-Text "54×25=" "62×57="

Replace-Text "27×84=" "72×86="
Replace-Text "59×63=" "70×13="
Replace-Text "33×73=" "48×91="
Replace-Text "14×45=" "12×82="
Replace-Text "82×67=" "49×23="

Replace-Text "18×93=" "74×54="
Replace-Text "15×69=" "49×11="
Replace-Text "68×74=" "75×73="
Replace-Text "98×25=" "24×93="
Replace-Text "46×19=" "99×94="

Replace-Text "59×79=" "80×43="
Replace-Text "25×60=" "83×80="
Replace-Text "46×12=" "85×37="
Replace-Text "55×44=" "75×87="
Replace-Text "95×42=" "71×54="
